$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new row 18 with the new item "PowerBank Xiaomi"
$ws.Range("A18").Value = 17
$ws.Range("B18").Value = "PowerBank Xiaomi"
$ws.Range("C18").Value = 600000
$ws.Range("D18").Value = 1000
$ws.Range("E18").Value = "iBox"
